$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# A new day's data was appended to the bottom of the table. The row that used
# to hold the trailing "統合" footnote (row 95) shifts down to row 96, so
# insert a fresh row at 95 (this also pushes the footnote row down and gives
# the new row the same formatting as the row above it, same as Excel's
# default Insert behaviour).
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new day's figures.
$ws.Cells.Item(95, 1).Value = 43950
$ws.Cells.Item(95, 2).Value = 396
$ws.Cells.Item(95, 3).Value = 31510
$ws.Cells.Item(95, 4).Value = 0
$ws.Cells.Item(95, 5).Value = 6664

# The sheet's print area grew by one row to keep including the footnote row.
foreach ($n in $wb.Names) {
  $n.RefersTo = "=" + $ws.Name + "!`$A`$1:`$E`$98"
}

# The active selection moves to the newly entered row.
[void]$ws.Range("A95").Select()
